$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-10-18 03:59:44"

$wsZhCn.Range("H3").Value = "2016-10-18 03:59:29"
$wsZhCn.Range("K3").Value = "2016-10-18 04:00:29"

$wsDeDe.Range("K3").Value = "2016-10-18 04:00:53"
